# Update Name of Algo
# Applies updated values to the result_data_RandomForest sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.3816
$ws.Range("C3").Value = -11.7091
$ws.Range("E4").Value = 12.62590000000001
$ws.Range("C5").Value = -12.01729999999999
$ws.Range("E6").Value = 11.98699999999999
$ws.Range("D7").Value = -6.858299999999993
$ws.Range("A9").Value = -20.40839999999998
$ws.Range("D9").Value = -8.719900000000003
$ws.Range("E10").Value = 11.60749999999999
$ws.Range("C11").Value = -13.80420000000001
$ws.Range("E11").Value = 13.2419
$ws.Range("C12").Value = -14.26030000000001
$ws.Range("A13").Value = -21.89230000000003
$ws.Range("A16").Value = -20.09299999999999
$ws.Range("A18").Value = -21.89920000000003
$ws.Range("A20").Value = -21.97140000000002
$ws.Range("C21").Value = -13.72380000000001
$ws.Range("D21").Value = -8.377100000000004
$ws.Range("E21").Value = 12.8099
$ws.Range("E25").Value = 13.2827
